# C5-PowerPoint.pptx edit
#
# 1. Slide 6 ("SOURCES OF FINANCE") table: switch the applied table style
#    from the default "Medium Style 2" ({0018A540-78F9-4BC9-94E2-1EE7A3C0060F})
#    to {6478EA24-861A-4F4C-9DE0-8418FC414E6A} (a different built-in table
#    style), the same way PowerPoint's Table Design > Table Styles gallery
#    would via Table.ApplyStyle.
#
# 2. The presentation's applied theme (design "Integral") is swapped back to
#    the stock "Office Theme" colour palette. Re-point every themed colour
#    slot (dk2/lt2/accent1-6/hlink/folHlink - dk1/lt1 are already identical)
#    back to the default Office values using ThemeColorScheme, which is the
#    supported, non-destructive way to edit theme colours in this host.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------
$tableSlide = $p.Slides.Item(6)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{6478EA24-861A-4F4C-9DE0-8418FC414E6A}")
    }
}

# --- 2. Theme colours -------------------------------------------------------
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
